# Further analysis of zooplankton blitz - crosswalk table cleanup on "zooper" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zooper")
$ws.Activate()

# 1. Remove the obsolete "Cumaceans" (CUMAC) crosswalk row. This was row 106;
#    deleting it shifts every subsequent row (107-119) up by one, and the
#    now-unreferenced "Cumacea" shared string is dropped automatically.
$ws.Rows.Item(106).Delete()

# 2. Fix a species-name typo: "Pseudodiaptomus forbesi" -> "Pseudodiaptomus forbesii"
$ws.Range("D31").Value = "Pseudodiaptomus forbesii"

# 3. Re-bucket several genus-level "Analy2" (column L) entries into their
#    coarser Order-level "Other" categories.
$ws.Range("L26").Value = "Calanoid Other"
$ws.Range("L80").Value = "Cyclopoid Other"
$ws.Range("L81").Value = "Cyclopoid Other"
$ws.Range("L82").Value = "Cyclopoid Other"
$ws.Range("L83").Value = "Cyclopoid Other"
$ws.Range("L98").Value = "Cyclopoid Other"
$ws.Range("L99").Value = "Cyclopoid Other"
$ws.Range("L100").Value = "Cyclopoid Other"

# 4. Widen column D (FRP_Meso) to fit the longer label text.
$ws.Columns.Item(4).ColumnWidth = 38.8

# 5. Leave the selection on the row that now occupies position 106
#    (previously row 107, the "crab zoea" entry), matching the editor's
#    last on-screen action.
$ws.Rows.Item(106).Select()
